# Add season-record columns (Wins / Losses / Ties) to the MIL_2000 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels in AD1:AF1, matching the existing
#     header formatting (bold font, centered/top alignment, thin border)
#     by copying the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-48): the team's season record, same for every player row.
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 89   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 1    # AF -> Ties
}
